# Follow change of xlsx.js 2.0.0.
#
# B3 used to hold a "serial-date-as-plain-number" value (41202, styled
# with the out-of-range/leftover style index 1) and is now written out as
# a proper date/time serial (41194.375) formatted with a real (built-in)
# date number format. Along the way two new cellXfs records come into
# existence: a currency-style record (numFmtId 5) and the date record
# actually used by B3 (numFmtId 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B3's current style index (s="1") points past the end of cellXfs in the
# source file, so touching it directly would blow up - snap it back to
# the built-in "Normal" style (index 0) first.
$ws.Range("B3").Style = "Normal"

# Create the (otherwise unused) numFmtId 5 ("$#,##0_);($#,##0)") cellXfs
# record via a scratch cell, then restore that cell exactly (value AND
# style) so only the style table grows - A1 itself must come out unchanged.
$scratch = $ws.Range("A1")
$scratchValue = $scratch.Value2
$scratch.NumberFormat = "$#,##0_);($#,##0)"
$scratch.Style = "Normal"
$scratch.Value = $scratchValue

# Give B3 the numFmtId 14 ("mm-dd-yy") date format and its new value.
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = 41194.375
